# Update the "scatter" worksheet's Y1 (column B) and Y2 (column C) values,
# rows 2-11, to match the new data set. The scatter chart on this sheet
# references these ranges directly, so its cached values will refresh
# automatically when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scatter")

$newB = @(7, 4, 7, 6, 8, 9, 1, 8, 0, 8)
$newC = @(6, 7, 5, 8, 3, 1, 2, 1, 5, 3)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newB[$i]
    $ws.Cells.Item($row, 3).Value = $newC[$i]
}
